$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Orders")

# Update text (shared string) values in column A (Title-ish ids) and column B (Title)
$ws.Range("A2").Value = "Duis culpa"
$ws.Range("A3").Value = "sint officia eiusmod nulla"
$ws.Range("A4").Value = "dolore laboris"

$ws.Range("B2").Value = "consectetur Ut voluptate et dolor"
$ws.Range("B3").Value = "dolore ad eiusmod"
$ws.Range("B4").Value = "voluptate"

# Update numeric values in columns C and D
$ws.Range("C2").Value = 294902898
$ws.Range("D2").Value = 690495123

$ws.Range("C3").Value = -942386301
$ws.Range("D3").Value = -1157565972

$ws.Range("C4").Value = -1992719263
$ws.Range("D4").Value = 982097709
